# Lab 6 / Lab 7 title updates
# - Slide 1 (Title 5): "Power BI Quick Start #3 " -> "Power BI Quick Start #2 ",
#                        "E02" -> "E03"
# - Slide 2 (Title 2): "E01" -> "E03" (leading "Power BI Quick Start #2 " run
#                        is left untouched)
#
# Character-range edits are used (rather than replacing the whole
# TextRange.Text) so each run keeps its own original run properties (rPr) -
# only the <a:t> text content changes, just like in the source commit.

$p = $ppt.ActivePresentation

# --- Slide 1: Title 5 -------------------------------------------------
$slide1 = $p.Slides.Item(1)
$title1 = $slide1.Shapes.Item(1)
$tr1 = $title1.TextFrame.TextRange

$run1a = $tr1.Characters(1, 24)   # "Power BI Quick Start #3 "
$run1b = $tr1.Characters(25, 3)   # "E02"

$run1a.Text = "Power BI Quick Start #2 "
$run1b.Text = "E03"

# --- Slide 2: Title 2 -------------------------------------------------
$slide2 = $p.Slides.Item(2)
$title2 = $slide2.Shapes.Item(1)
$tr2 = $title2.TextFrame.TextRange

$run2b = $tr2.Characters(25, 3)   # "E01"
$run2b.Text = "E03"
